$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Insert two new rows at position 4-5 (pushing the existing "dataset.commit.*"
# rows, and everything below, down by two) to make room for the new
# dataset.preview.table / dataset.preview.line entries.
$ws.Range("A4:B5").EntireRow.Insert(-4121, $null)

$tableFormula = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"
$lineFormula  = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

$ws.Cells.Item(4, 1).Value2 = "dataset.preview.table"
$ws.Cells.Item(4, 2).Value2 = $tableFormula
$ws.Cells.Item(5, 1).Value2 = "dataset.preview.line"
$ws.Cells.Item(5, 2).Value2 = $lineFormula

# New rows get wrapped text, vertically centered, with a tall 120pt row
# height so the multi-line preview formulas are fully visible.
$previewRange = $ws.Range("A4:B5")
$previewRange.WrapText = $true
$previewRange.VerticalAlignment = -4108

$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 120

# Scroll/selection ends up on B7 (was B35 before the two rows were inserted
# above it), and the old frozen topLeftCell="A13" view is gone.
$ws.Range("B7").Select() | Out-Null
